# Fix: set server.js as entry point and add start script
# (adds the two new feedback rows captured after the fix)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feedbacks")

# Plain text values (Excel already keeps these as text by default).
function Set-Text($addr, $text) {
    $ws.Range($addr).Value = $text
}

# Values that look like pure numbers ("2020") must be forced to stay
# text - otherwise Excel auto-converts them to a numeric cell. Flip the
# cell to Text format just long enough to take the string, then drop the
# formatting back to Normal so the cell keeps the default (unstyled) xf.
function Set-TextNumeric($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 14
Set-Text        "A14" "20BEE2001"
Set-Text        "B14" "Ram"
Set-TextNumeric "C14" "2020"
Set-Text        "D14" "Deis Irae"
$ws.Range("E14").Value = 5
Set-Text        "F14" "5/2/2026, 11:25:50 am"

# Row 15
Set-Text        "A15" "20BEE2001"
Set-Text        "B15" "Norman"
Set-TextNumeric "C15" "2020"
Set-Text        "D15" "Canteen is really clean"
$ws.Range("E15").Value = 5
Set-Text        "F15" "5/2/2026, 11:47:10 am"
